# gh-pages data refresh: re-scrape of bilibili exhibition listings (output generated at 456a3b4).
# Rows 2-10 on the "展览" and "全部类型" sheets are refreshed in place with the newest 9 events.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
  $ws = $wb.Worksheets.Item($sheetName)

  $cB = $ws.Cells.Item(2, 2)
  $cB.NumberFormat = "@"
  $cB.Value = "2024-07-20"
  $cB.Style = "Normal"
  $ws.Cells.Item(2, 3).Value = "丽水·CCAC动漫游戏嘉年华"
  $ws.Cells.Item(2, 4).Value = "南环西路109号 九城宴会中心"
  $ws.Cells.Item(2, 5).Value = "2024.07.20 09:00-07.20 16:00"
  $ws.Cells.Item(2, 6).Value = 123
  $ws.Cells.Item(2, 7).Value = 50
  $ws.Cells.Item(2, 8).Value = "https://show.bilibili.com/platform/detail.html?id=86306"
  $ws.Cells.Item(2, 9).Value = "//i1.hdslb.com/bfs/openplatform/202405/4TrBjBlV1716551375116.png"

  $cB = $ws.Cells.Item(3, 2)
  $cB.NumberFormat = "@"
  $cB.Value = "2024-07-20"
  $cB.Style = "Normal"
  $ws.Cells.Item(3, 3).Value = "丽水·龙泉ACG动漫游戏博览会"
  $ws.Cells.Item(3, 4).Value = "南秦路1号望瓯·陶溪川直走200米左手边(7号楼) 望瓯陶溪川活动中心"
  $ws.Cells.Item(3, 5).Value = "2024.07.20 10:00-07.21 18:00"
  $ws.Cells.Item(3, 6).Value = 1650
  $ws.Cells.Item(3, 7).Value = 55
  $ws.Cells.Item(3, 8).Value = "https://show.bilibili.com/platform/detail.html?id=86671"
  $ws.Cells.Item(3, 9).Value = "//i0.hdslb.com/bfs/openplatform/202406/LSorIT7S1717486817969.png"

  $cB = $ws.Cells.Item(4, 2)
  $cB.NumberFormat = "@"
  $cB.Value = "2024-07-20"
  $cB.Style = "Normal"
  $ws.Cells.Item(4, 3).Value = "青田·苍渊动漫游戏嘉年华"
  $ws.Cells.Item(4, 4).Value = "鹤城街道高湾1号 青田侨乡世茂大酒店"
  $ws.Cells.Item(4, 5).Value = "2024.07.20 10:00-07.20 17:00"
  $ws.Cells.Item(4, 6).Value = 21
  $ws.Cells.Item(4, 7).Value = 45
  $ws.Cells.Item(4, 8).Value = "https://show.bilibili.com/platform/detail.html?id=88501"
  $ws.Cells.Item(4, 9).Value = "//i2.hdslb.com/bfs/openplatform/202406/fLPkC9eX1718813133406.jpeg"

  $cB = $ws.Cells.Item(5, 2)
  $cB.NumberFormat = "@"
  $cB.Value = "2024-07-27"
  $cB.Style = "Normal"
  $ws.Cells.Item(5, 3).Value = "丽水·thp01～风摄少微"
  $ws.Cells.Item(5, 4).Value = "大猷街 应星楼"
  $ws.Cells.Item(5, 5).Value = "2024.07.27 10:00-07.27 18:00"
  $ws.Cells.Item(5, 6).Value = 24
  $ws.Cells.Item(5, 7).Value = 50
  $ws.Cells.Item(5, 8).Value = "https://show.bilibili.com/platform/detail.html?id=87134"
  $ws.Cells.Item(5, 9).Value = "//i2.hdslb.com/bfs/openplatform/202406/JuvSmncN1717775885615.png"

  $cB = $ws.Cells.Item(6, 2)
  $cB.NumberFormat = "@"
  $cB.Value = "2024-07-27"
  $cB.Style = "Normal"
  $ws.Cells.Item(6, 3).Value = "丽水·第四届HP国风动漫游戏嘉年华"
  $ws.Cells.Item(6, 4).Value = "城北街798号 莱茵体育生活馆"
  $ws.Cells.Item(6, 5).Value = "2024.07.27 08:30-07.27 17:00"
  $ws.Cells.Item(6, 6).Value = 427
  $ws.Cells.Item(6, 7).Value = 65
  $ws.Cells.Item(6, 8).Value = "https://show.bilibili.com/platform/detail.html?id=87305"
  $ws.Cells.Item(6, 9).Value = "//i2.hdslb.com/bfs/openplatform/202406/YUnPOKGV1718268952725.jpeg"

  $cB = $ws.Cells.Item(7, 2)
  $cB.NumberFormat = "@"
  $cB.Value = "2024-08-03"
  $cB.Style = "Normal"
  $ws.Cells.Item(7, 3).Value = "丽水·樱卡动漫游戏嘉年华"
  $ws.Cells.Item(7, 4).Value = "中东路848号(解放街交汇) 飞达国际大酒店"
  $ws.Cells.Item(7, 5).Value = "2024.08.03 10:00-08.03 17:00"
  $ws.Cells.Item(7, 6).Value = 148
  $ws.Cells.Item(7, 7).Value = 50
  $ws.Cells.Item(7, 8).Value = "https://show.bilibili.com/platform/detail.html?id=87276"
  $ws.Cells.Item(7, 9).Value = "//i0.hdslb.com/bfs/openplatform/202406/bVp0Zg1B1718172430380.jpeg"

  $cB = $ws.Cells.Item(8, 2)
  $cB.NumberFormat = "@"
  $cB.Value = "2024-08-10"
  $cB.Style = "Normal"
  $ws.Cells.Item(8, 3).Value = "丽水·CCAC动漫七夕（回馈展）"
  $ws.Cells.Item(8, 4).Value = "中东路848号(解放街交汇) 飞达国际大酒店"
  $ws.Cells.Item(8, 5).Value = "2024.08.10 09:00-08.10 17:00"
  $ws.Cells.Item(8, 6).Value = 67
  $ws.Cells.Item(8, 7).Value = 29.9
  $ws.Cells.Item(8, 8).Value = "https://show.bilibili.com/platform/detail.html?id=86567"
  $ws.Cells.Item(8, 9).Value = "//i0.hdslb.com/bfs/openplatform/202405/tsOzbBRx1717015539538.png"

  $cB = $ws.Cells.Item(9, 2)
  $cB.NumberFormat = "@"
  $cB.Value = "2024-08-17"
  $cB.Style = "Normal"
  $ws.Cells.Item(9, 3).Value = "丽水·AEO纯白礼赞动漫嘉年华"
  $ws.Cells.Item(9, 4).Value = "城北街1001号 爱依·时尚婚宴中心"
  $ws.Cells.Item(9, 5).Value = "2024.08.17 09:00-08.17 16:00"
  $ws.Cells.Item(9, 6).Value = 532
  $ws.Cells.Item(9, 7).Value = 55
  $ws.Cells.Item(9, 8).Value = "https://show.bilibili.com/platform/detail.html?id=86779"
  $ws.Cells.Item(9, 9).Value = "//i2.hdslb.com/bfs/openplatform/202406/MxJ3oNjt1717405405850.jpeg"

  $cB = $ws.Cells.Item(10, 2)
  $cB.NumberFormat = "@"
  $cB.Value = "2024-09-16"
  $cB.Style = "Normal"
  $ws.Cells.Item(10, 3).Value = "丽水·LZ栗子动漫游戏嘉年华"
  $ws.Cells.Item(10, 4).Value = "城北街798号 莱茵体育生活馆"
  $ws.Cells.Item(10, 5).Value = "2024.09.16 09:30-09.16 17:00"
  $ws.Cells.Item(10, 6).Value = 399
  $ws.Cells.Item(10, 7).Value = 65
  $ws.Cells.Item(10, 8).Value = "https://show.bilibili.com/platform/detail.html?id=87480"
  $ws.Cells.Item(10, 9).Value = "//i1.hdslb.com/bfs/openplatform/202406/bATqcZhH1719285865931.jpeg"

}

